$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - first worksheet
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 71
$wsExhibit.Range("F4").Value = 1505
$wsExhibit.Range("F5").Value = 579
$wsExhibit.Range("F6").Value = 1068
$wsExhibit.Range("F7").Value = 11062
$wsExhibit.Range("F8").Value = 11062
$wsExhibit.Range("F10").Value = 32
$wsExhibit.Range("F11").Value = 320
$wsExhibit.Range("F12").Value = 1066
$wsExhibit.Range("F13").Value = 754
$wsExhibit.Range("F14").Value = 12225
$wsExhibit.Range("F15").Value = 12763
$wsExhibit.Range("F22").Value = 39

# Sheet "全部类型" (All types) - fourth worksheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 71
$wsAll.Range("F5").Value = 1505
$wsAll.Range("F6").Value = 579
$wsAll.Range("F7").Value = 1068
$wsAll.Range("F8").Value = 11062
$wsAll.Range("F9").Value = 11062
$wsAll.Range("F11").Value = 32
$wsAll.Range("F12").Value = 320
$wsAll.Range("F13").Value = 1066
$wsAll.Range("F14").Value = 754
$wsAll.Range("F15").Value = 12225
$wsAll.Range("F16").Value = 12763
$wsAll.Range("F23").Value = 39
